$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the "Verb" header in A1 to "#Verb"
$ws.Range("A1").Value = "#Verb"

# Reset selection to the default (A1) instead of the previously selected E13
$ws.Range("A1").Select()
